# Updated lifelines in the two diagrams
#
# 1) Extend the "Model" object's lifeline with a new dashed connector
#    underneath it (new "Straight Connector 57" shape).
# 2) Extend the existing left-most lifeline ("Straight Connector 4") so it
#    reaches further down the sequence diagram.
# 3) Nudge the "execute(\u201creverse\u201d)" call-out textbox a bit further left.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- 1) Add the new lifeline connector under the "Model" box -----------
# Duplicate an existing lifeline connector so the new shape inherits the
# exact same line style / theme style refs (color, dash, weight, style).
$template = $s.Shapes.Item("Straight Connector 4")
$newCxn = $template.Duplicate().Item(1)
$newCxn.Name = "Straight Connector 57"

# Position / size (EMU 5332155,4597186 / 0 x 1073583), expressed in points
# with the extra digits needed to survive the COM Single round-trip.
$newCxn.Left = 419.854736328125
$newCxn.Top = 361.983154296875
$newCxn.Width = 0.0
$newCxn.Height = 84.53409576416016

# Send it to the back of the z-order so it lands right after the group
# shape properties, before "Rectangle 65" - matching the original author's
# draw order.
$newCxn.ZOrder(1)

# --- 2) Stretch the existing lifeline further down ----------------------
$lifeline = $s.Shapes.Item("Straight Connector 4")
$lifeline.Height = 397.8401794433594

# --- 3) Shift the call-out textbox left ---------------------------------
$callout = $s.Shapes.Item("TextBox 25")
$callout.Left = -78.0
